$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps its original text representation
# (values such as "1.00" or "0.230" would otherwise be auto-converted
# to numbers by Excel, losing the trailing zeros / exact formatting).
$ws.Range("D2:D51").NumberFormat = "@"

# Price (column D) updates
$ws.Range("D2").Value = "91.406.06"
$ws.Range("D3").Value = "3.154.06"
$ws.Range("D5").Value = "238.91"
$ws.Range("D6").Value = "621.10"
$ws.Range("D7").Value = "1.12"
$ws.Range("D8").Value = "0.375"
$ws.Range("D9").Value = "1.00"
$ws.Range("D10").Value = "3.155.14"
$ws.Range("D11").Value = "0.744"
$ws.Range("D13").Value = "0.0000246"
$ws.Range("D14").Value = "35.45"
$ws.Range("D16").Value = "91.164.81"
$ws.Range("D17").Value = "3.750.65"
$ws.Range("D18").Value = "3.177.24"
$ws.Range("D19").Value = "3.75"
$ws.Range("D20").Value = "15.30"
$ws.Range("D22").Value = "457.22"
$ws.Range("D23").Value = "0.0000204"
$ws.Range("D24").Value = "9.18"
$ws.Range("D25").Value = "6.03"
$ws.Range("D26").Value = "89.18"
$ws.Range("D27").Value = "12.05"
$ws.Range("D29").Value = "0.997"
$ws.Range("D30").Value = "0.128"
$ws.Range("D32").Value = "0.230"
$ws.Range("D33").Value = "9.38"
$ws.Range("D35").Value = "0.932"
$ws.Range("D36").Value = "7.67"
$ws.Range("D37").Value = "26.53"
$ws.Range("D38").Value = "513.19"
$ws.Range("D39").Value = "1.95"
$ws.Range("D40").Value = "1.35"
$ws.Range("D42").Value = "3.84"
$ws.Range("D43").Value = "3.46"
$ws.Range("D44").Value = "22.19"
$ws.Range("D46").Value = "0.712"
$ws.Range("D47").Value = "157.89"
$ws.Range("D49").Value = "1.37"
$ws.Range("D50").Value = "4.49"
$ws.Range("D51").Value = "44.08"

# Volume(1h) (column E) updates
$ws.Range("E2").Value = "  +2.23%  "
$ws.Range("E3").Value = "  +2.58%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("E5").Value = "  +1.70%  "
$ws.Range("E6").Value = "  +0.90%  "
$ws.Range("E7").Value = "  +6.20%  "
$ws.Range("E8").Value = "  +4.01%  "
$ws.Range("E9").Value = "  -0.06%  "
$ws.Range("E10").Value = "  +2.63%  "
$ws.Range("E11").Value = "  +5.13%  "
$ws.Range("E12").Value = "  +2.49%  "
$ws.Range("E13").Value = "  -1.24%  "
$ws.Range("E14").Value = "  +1.30%  "
$ws.Range("E15").Value = "  +4.57%  "
$ws.Range("E16").Value = "  +2.04%  "
$ws.Range("E17").Value = "  +2.93%  "
$ws.Range("E18").Value = "  +2.99%  "
$ws.Range("E19").Value = "  -0.21%  "
$ws.Range("E20").Value = "  +11.64%  "
$ws.Range("E21").Value = "  +11.13%  "
$ws.Range("E22").Value = "  +6.37%  "
$ws.Range("E23").Value = "  -3.02%  "
$ws.Range("E24").Value = "  +5.54%  "
$ws.Range("E25").Value = "  +5.66%  "
$ws.Range("E26").Value = "  +3.06%  "
$ws.Range("E27").Value = "  +3.45%  "
$ws.Range("E29").Value = "  -0.35%  "
$ws.Range("E30").Value = "  +43.92%  "
$ws.Range("E31").Value = "  +10.37%  "
$ws.Range("E32").Value = "  +17.94%  "
$ws.Range("E33").Value = "  +4.34%  "
$ws.Range("E34").Value = "  +13.85%  "
$ws.Range("E35").Value = "  -19.22%  "
$ws.Range("E36").Value = "  +8.37%  "
$ws.Range("E37").Value = "  +4.06%  "
$ws.Range("E38").Value = "  +4.46%  "
$ws.Range("E39").Value = "  +4.19%  "
$ws.Range("E40").Value = "  +8.49%  "
$ws.Range("E41").Value = "  +14.71%  "
$ws.Range("E42").Value = "  +5.58%  "
$ws.Range("E43").Value = "  -4.45%  "
$ws.Range("E44").Value = "  +0.58%  "
$ws.Range("E45").Value = "  -0.08%  "
$ws.Range("E46").Value = "  +6.72%  "
$ws.Range("E47").Value = "  +3.38%  "
$ws.Range("E48").Value = "  +5.36%  "
$ws.Range("E49").Value = "  +6.33%  "
$ws.Range("E50").Value = "  +4.91%  "
